$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from the last existing data row (53) down to the new row (54)
$ws.Range("A53:E53").Copy($ws.Range("A54:E54"))

# Populate the new row's values
$ws.Range("A54").Value = 45986
$ws.Range("B54").Value = 2025
$ws.Range("C54").Value = 2.043309689777173
$ws.Range("D54").Value = 2026
$ws.Range("E54").Value = 1.002299702378884
